# Sort the comma-separated "Recorded By" values in column G alphabetically
# using strict ASCII/ordinal ordering (matches the upstream sync edit,
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -gt 1) {
        $list = New-Object System.Collections.Generic.List[string]
        foreach ($p in $parts) { [void]$list.Add($p) }
        $list.Sort([System.StringComparer]::Ordinal)
        $newVal = [string]::Join(", ", $list)

        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
